# Add new columns I ("I0") and J ("IF") to the sheet, matching the
# formatting of the existing header row (copied from H1, which carries
# the bold/centered/bordered header style), and fill in the per-row
# data values for rows 2-77.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 onto the two new header cells so
# they pick up the same cell style (bold font, border, centered/top
# alignment) instead of Excel's default style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for column I (I0) and column J (IF), rows 2-77.
$iVals = @(7,7,7,8,7,8,7,6,7,8,8,8,8,8,7,6,8,8,8,8,8,9,9,9,9,8,8,8,8,8,8,8,8,8,8,8,7,8,8,8,8,8,8,8,8,7,7,7,7,8,6,5,7,8,8,7,8,8,8,8,8,7,6,9,8,8,8,6,6,7,6,8,8,6,6,4)
$jVals = @(7,7,7,8,7,9,7,6,7,8,8,9,8,8,7,7,8,8,8,8,8,9,9,9,9,8,8,8,8,8,8,8,8,9,8,8,8,8,8,8,8,8,8,8,9,8,8,8,8,8,6,6,7,8,8,8,8,9,8,8,8,8,7,9,8,8,8,6,6,7,7,8,8,6,6,4)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
